$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 15 - this pushes the existing rows 15..42 down to 16..43,
# carrying all of their data (and the D-column date style) with them.
$ws.Rows(15).Insert()

# Populate the newly inserted row 15 with the new weekly record.
$ws.Range("A15").Value = 10
$ws.Range("B15").Value = "Vega Modelo de Temuco"
$ws.Range("C15").Value = "La Araucanía"
$ws.Range("D15").Value = 44495
$ws.Range("E15").Value = 9
$ws.Range("F15").Value = 300000000
$ws.Range("G15").Value = "Espárragos"
$ws.Range("H15").Value = "Sin especificar"
$ws.Range("I15").Value = "Primera"
$ws.Range("J15").Value = 200
$ws.Range("K15").Value = 1300
$ws.Range("L15").Value = 1300
$ws.Range("M15").Value = 1300
$ws.Range("N15").Value = "$/kilo"
$ws.Range("O15").Value = "Región del Maule"
$ws.Range("P15").Value = 1300
$ws.Range("Q15").Value = 1
$ws.Range("R15").Value = "Hortaliza"
